# Update the "想去人数" (want-to-go count) values in column F
# for both the "展览" sheet and the "全部类型" sheet.
# These two sheets list overlapping events, so the same updates are
# applied to each sheet independently.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (rows 2-13, row 8 untouched) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value  = 601
$wsExhibit.Range("F3").Value  = 203
$wsExhibit.Range("F4").Value  = 474
$wsExhibit.Range("F5").Value  = 495
$wsExhibit.Range("F6").Value  = 282
$wsExhibit.Range("F7").Value  = 2554
$wsExhibit.Range("F9").Value  = 6901
$wsExhibit.Range("F10").Value = 186
$wsExhibit.Range("F11").Value = 436
$wsExhibit.Range("F12").Value = 6
$wsExhibit.Range("F13").Value = 34

# --- Sheet "全部类型" (rows 2-17, rows 7,8,10,15,16 untouched) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 601
$wsAll.Range("F3").Value  = 203
$wsAll.Range("F4").Value  = 474
$wsAll.Range("F5").Value  = 495
$wsAll.Range("F6").Value  = 282
$wsAll.Range("F9").Value  = 2554
$wsAll.Range("F11").Value = 6901
$wsAll.Range("F12").Value = 186
$wsAll.Range("F13").Value = 436
$wsAll.Range("F14").Value = 6
$wsAll.Range("F17").Value = 34
